$wb = $excel.ActiveWorkbook

# Delete the "Info" worksheet
$excel.DisplayAlerts = $false
$infoSheet = $wb.Worksheets.Item("Info")
$infoSheet.Delete()
$excel.DisplayAlerts = $true

# Rename the "scenario" worksheet to "Sheet1"
$scenarioSheet = $wb.Worksheets.Item("scenario")
$scenarioSheet.Name = "Sheet1"
